# Insert two new weekly price rows (Femacal de La Calera, Ají) at the top of
# the data block (rows 171-172), pushing the existing rows 171..262 down to
# 173..264.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 171 (shifts 171..262 down to 173..264).
$ws.Range("A171:R172").Insert()

# --- New row 171 ---
$ws.Cells.Item(171, 1).Value = 3
$ws.Cells.Item(171, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44452
$ws.Cells.Item(171, 5).Value = 5
$ws.Cells.Item(171, 6).Value = 100112021
$ws.Cells.Item(171, 7).Value = "Ají"
$ws.Cells.Item(171, 8).Value = "Americana (o)"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 45
$ws.Cells.Item(171, 11).Value = 72000
$ws.Cells.Item(171, 12).Value = 73000
$ws.Cells.Item(171, 13).Value = 72556
$ws.Cells.Item(171, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(171, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(171, 16).Value = 2902
$ws.Cells.Item(171, 17).Value = 25
$ws.Cells.Item(171, 18).Value = "Hortaliza"

# --- New row 172 ---
$ws.Cells.Item(172, 1).Value = 3
$ws.Cells.Item(172, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44452
$ws.Cells.Item(172, 5).Value = 5
$ws.Cells.Item(172, 6).Value = 100112021
$ws.Cells.Item(172, 7).Value = "Ají"
$ws.Cells.Item(172, 8).Value = "Inferno"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 73
$ws.Cells.Item(172, 11).Value = 41000
$ws.Cells.Item(172, 12).Value = 42000
$ws.Cells.Item(172, 13).Value = 41521
$ws.Cells.Item(172, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 2768
$ws.Cells.Item(172, 17).Value = 15
$ws.Cells.Item(172, 18).Value = "Hortaliza"
